# Updated symbol list on Tue Feb  7 07:42:49 UTC 2023 with GitHub Actions
# Refresh the Price (column D) and Volume(1h) (column E) values for the
# crypto rows whose figures moved since the last scrape. Values are stored
# as literal text (inlineStr) in the sheet, so each cell is forced to the
# Text number format before the assignment to stop Excel from re-parsing
# numeric- / percent-looking strings back into real numbers (which would
# also silently drop significant trailing zeros, e.g. "0.0002000").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "328.68"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.09%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "43.97"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.28%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.518"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.53%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08008"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.62%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.989"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.93%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.27%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.574"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-4.81%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9498"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.02%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1125"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-3.92%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.83%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "10.66"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "26.52%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09981"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.05%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04781"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "13.48%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.07%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001271"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.87%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-3.94%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.005980"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.70%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.363"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-6.20%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3474"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.69%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1420"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "3.06%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.74%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.46%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004324"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-5.20%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001201"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.80%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003745"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-6.10%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02597"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-1.49%"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.02%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007555"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.58%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1396"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.04%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007407"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "3.22%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002016"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.69%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008634"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-5.99%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.08%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.06%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.003531"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "55.57%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003785"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "7.35%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.06%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.06%"
